# Update the heading date, then each division-answer cell in the table.
# Find/Replace is applied in document order so that the one coincidental
# collision (an original answer equals a later replacement's new text,
# "39÷4=9, 3") resolves correctly: the cell holding that text is changed
# to its own new value before the later cell is searched/replaced to the
# same string.
$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-01-07 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-01-08 Monday", 2) | Out-Null
$d.Content.Find.Execute("49÷4=12, 1", $true, $false, $false, $false, $false, $true, 1, $false, "36÷9=4, 0", 2) | Out-Null
$d.Content.Find.Execute("80÷5=16, 0", $true, $false, $false, $false, $false, $true, 1, $false, "30÷2=15, 0", 2) | Out-Null
$d.Content.Find.Execute("63÷5=12, 3", $true, $false, $false, $false, $false, $true, 1, $false, "72÷3=24, 0", 2) | Out-Null
$d.Content.Find.Execute("34÷8=4, 2", $true, $false, $false, $false, $false, $true, 1, $false, "23÷9=2, 5", 2) | Out-Null
$d.Content.Find.Execute("70÷5=14, 0", $true, $false, $false, $false, $false, $true, 1, $false, "24÷2=12, 0", 2) | Out-Null
$d.Content.Find.Execute("39÷4=9, 3", $true, $false, $false, $false, $false, $true, 1, $false, "88÷5=17, 3", 2) | Out-Null
$d.Content.Find.Execute("24÷5=4, 4", $true, $false, $false, $false, $false, $true, 1, $false, "33÷3=11, 0", 2) | Out-Null
$d.Content.Find.Execute("18÷6=3, 0", $true, $false, $false, $false, $false, $true, 1, $false, "70÷9=7, 7", 2) | Out-Null
$d.Content.Find.Execute("42÷4=10, 2", $true, $false, $false, $false, $false, $true, 1, $false, "84÷5=16, 4", 2) | Out-Null
$d.Content.Find.Execute("77÷3=25, 2", $true, $false, $false, $false, $false, $true, 1, $false, "32÷9=3, 5", 2) | Out-Null
$d.Content.Find.Execute("17÷2=8, 1", $true, $false, $false, $false, $false, $true, 1, $false, "66÷7=9, 3", 2) | Out-Null
$d.Content.Find.Execute("17÷5=3, 2", $true, $false, $false, $false, $false, $true, 1, $false, "39÷4=9, 3", 2) | Out-Null
$d.Content.Find.Execute("85÷6=14, 1", $true, $false, $false, $false, $false, $true, 1, $false, "87÷6=14, 3", 2) | Out-Null
$d.Content.Find.Execute("77÷7=11, 0", $true, $false, $false, $false, $false, $true, 1, $false, "90÷5=18, 0", 2) | Out-Null
$d.Content.Find.Execute("14÷9=1, 5", $true, $false, $false, $false, $false, $true, 1, $false, "53÷8=6, 5", 2) | Out-Null
$d.Content.Find.Execute("97÷8=12, 1", $true, $false, $false, $false, $false, $true, 1, $false, "85÷4=21, 1", 2) | Out-Null
$d.Content.Find.Execute("67÷8=8, 3", $true, $false, $false, $false, $false, $true, 1, $false, "93÷7=13, 2", 2) | Out-Null
$d.Content.Find.Execute("26÷3=8, 2", $true, $false, $false, $false, $false, $true, 1, $false, "56÷8=7, 0", 2) | Out-Null
$d.Content.Find.Execute("56÷4=14, 0", $true, $false, $false, $false, $false, $true, 1, $false, "71÷4=17, 3", 2) | Out-Null
$d.Content.Find.Execute("11÷4=2, 3", $true, $false, $false, $false, $false, $true, 1, $false, "35÷7=5, 0", 2) | Out-Null
$d.Content.Find.Execute("55÷2=27, 1", $true, $false, $false, $false, $false, $true, 1, $false, "44÷5=8, 4", 2) | Out-Null
$d.Content.Find.Execute("11÷6=1, 5", $true, $false, $false, $false, $false, $true, 1, $false, "80÷7=11, 3", 2) | Out-Null
$d.Content.Find.Execute("20÷8=2, 4", $true, $false, $false, $false, $false, $true, 1, $false, "29÷3=9, 2", 2) | Out-Null
$d.Content.Find.Execute("67÷2=33, 1", $true, $false, $false, $false, $false, $true, 1, $false, "94÷9=10, 4", 2) | Out-Null
$d.Content.Find.Execute("79÷5=15, 4", $true, $false, $false, $false, $false, $true, 1, $false, "79÷3=26, 1", 2) | Out-Null
